$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "code': 139, 'message': 'You have already favorited this status.'"
$ws.Range("A4").Select()
